$d = $word.ActiveDocument

# Locate the paragraph that contains the "m:'\n'" field (a Word field built
# from fldChar begin/instrText.../fldChar end), the one the M2Doc query
# template uses for the newLine() call. We search paragraph-by-paragraph so
# the script does not depend on a hard-coded field index.
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    foreach ($f in $p.Range.Fields) {
        if ($f.Code.Text -match "m:'\\n'") {
            $targetPara = $p
        }
    }
}

if ($targetPara -ne $null) {
    # Range covering the whole field (fldChar begin ... fldChar end) but not
    # the trailing paragraph mark, so the paragraph itself is preserved.
    $r = $d.Range($targetPara.Range.Start, $targetPara.Range.End - 1)

    # Replace the field with 5 plain-text runs spelling out the same
    # characters the field code used to carry ("{", "m", ":'\n", "'", "}"),
    # turning the M2Doc field code into literal template text, as done by
    # TokenIteratorFieldRewriterSplit.
    $xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p w:rsidP="00F5495F" w:rsidR="00C52979" w:rsidRDefault="00C52979"><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:'\n</w:t></w:r><w:r><w:t>'</w:t></w:r><w:r><w:t xml:space="preserve">}</w:t></w:r></w:p>
</w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>
"@
    $r.InsertXML($xml)
}
